# Final cleanup for the 2E Python section: drop the unused "_1" column
# (old column C, which only ever held a placeholder header and no data)
# and rename the two stat-value headers so they read "Australia_priceprice"
# and "Australia_pointspoints".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the empty placeholder column; this shifts the old "points" column
# (D) left into C, carrying its data and formatting with it.
$ws.Columns("C").Delete()

# Rename headers to match the cleaned-up naming.
$ws.Range("B1").Value = "Australia_priceprice"
$ws.Range("C1").Value = "Australia_pointspoints"

# Re-apply the intended column widths now that there are only 3 columns.
$ws.Columns("B").ColumnWidth = 21.83
$ws.Columns("C").ColumnWidth = 23.83
